$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = [DateTime]"2020-12-04"
$ws.Cells.Item(2, 11).Value = "Castle Brite"
$ws.Cells.Item(2, 12).Value = "Segunda"
$ws.Cells.Item(2, 13).Value = 500
$ws.Cells.Item(2, 14).Value = 15000
$ws.Cells.Item(2, 15).Value = 16000
$ws.Cells.Item(2, 16).Value = 15500
$ws.Cells.Item(2, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(2, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(2, 19).Value = 1033
$ws.Cells.Item(2, 20).Value = 15

# Row 3
$ws.Cells.Item(3, 4).Value = [DateTime]"2020-12-29"
$ws.Cells.Item(3, 11).Value = "Castle Brite"
$ws.Cells.Item(3, 12).Value = "Segunda"
$ws.Cells.Item(3, 13).Value = 300
$ws.Cells.Item(3, 14).Value = 15000
$ws.Cells.Item(3, 15).Value = 16000
$ws.Cells.Item(3, 16).Value = 15500
$ws.Cells.Item(3, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(3, 18).Value = "Región Metropolitana"
$ws.Cells.Item(3, 19).Value = 1033
$ws.Cells.Item(3, 20).Value = 15

# Row 4
$ws.Cells.Item(4, 4).Value = [DateTime]"2022-01-14"
$ws.Cells.Item(4, 11).Value = "Modesto"
$ws.Cells.Item(4, 12).Value = "Especial"
$ws.Cells.Item(4, 13).Value = 200
$ws.Cells.Item(4, 14).Value = 21000
$ws.Cells.Item(4, 15).Value = 21000
$ws.Cells.Item(4, 16).Value = 21000
$ws.Cells.Item(4, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(4, 18).Value = "Región Metropolitana"
$ws.Cells.Item(4, 19).Value = 1167
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value = [DateTime]"2022-01-14"
$ws.Cells.Item(5, 11).Value = "Modesto"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 18000
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 18000
$ws.Cells.Item(5, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(5, 18).Value = "Región Metropolitana"
$ws.Cells.Item(5, 19).Value = 1000
$ws.Cells.Item(5, 20).Value = 18

# Row 6
$ws.Cells.Item(6, 4).Value = [DateTime]"2022-01-14"
$ws.Cells.Item(6, 11).Value = "Modesto"
$ws.Cells.Item(6, 12).Value = "Segunda"
$ws.Cells.Item(6, 13).Value = 200
$ws.Cells.Item(6, 14).Value = 16000
$ws.Cells.Item(6, 15).Value = 16000
$ws.Cells.Item(6, 16).Value = 16000
$ws.Cells.Item(6, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(6, 18).Value = "Región Metropolitana"
$ws.Cells.Item(6, 19).Value = 889
$ws.Cells.Item(6, 20).Value = 18

# Row 7
$ws.Cells.Item(7, 4).Value = [DateTime]"2021-12-17"
$ws.Cells.Item(7, 11).Value = "Castle Brite"
$ws.Cells.Item(7, 12).Value = "Especial"
$ws.Cells.Item(7, 13).Value = 350
$ws.Cells.Item(7, 14).Value = 20000
$ws.Cells.Item(7, 15).Value = 20000
$ws.Cells.Item(7, 16).Value = 20000
$ws.Cells.Item(7, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(7, 18).Value = "Región Metropolitana"
$ws.Cells.Item(7, 19).Value = 1111
$ws.Cells.Item(7, 20).Value = 18

# Row 8
$ws.Cells.Item(8, 4).Value = [DateTime]"2021-12-17"
$ws.Cells.Item(8, 11).Value = "Castle Brite"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 350
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(8, 18).Value = "Región Metropolitana"
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 18

# Row 9
$ws.Cells.Item(9, 4).Value = [DateTime]"2021-12-17"
$ws.Cells.Item(9, 11).Value = "Castle Brite"
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 350
$ws.Cells.Item(9, 14).Value = 16000
$ws.Cells.Item(9, 15).Value = 16000
$ws.Cells.Item(9, 16).Value = 16000
$ws.Cells.Item(9, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(9, 18).Value = "Región Metropolitana"
$ws.Cells.Item(9, 19).Value = 889
$ws.Cells.Item(9, 20).Value = 18

# Row 10
$ws.Cells.Item(10, 4).Value = [DateTime]"2020-12-11"
$ws.Cells.Item(10, 11).Value = "Castle Brite"
$ws.Cells.Item(10, 12).Value = "Segunda"
$ws.Cells.Item(10, 13).Value = 500
$ws.Cells.Item(10, 14).Value = 15000
$ws.Cells.Item(10, 15).Value = 16000
$ws.Cells.Item(10, 16).Value = 15500
$ws.Cells.Item(10, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(10, 18).Value = "Región Metropolitana"
$ws.Cells.Item(10, 19).Value = 1033
$ws.Cells.Item(10, 20).Value = 15

# Row 11
$ws.Cells.Item(11, 4).Value = [DateTime]"2020-11-24"
$ws.Cells.Item(11, 11).Value = "Castle Brite"
$ws.Cells.Item(11, 12).Value = "Tercera"
$ws.Cells.Item(11, 13).Value = 400
$ws.Cells.Item(11, 14).Value = 15500
$ws.Cells.Item(11, 15).Value = 16000
$ws.Cells.Item(11, 16).Value = 15750
$ws.Cells.Item(11, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 1050
$ws.Cells.Item(11, 20).Value = 15

# Row 12
$ws.Cells.Item(12, 4).Value = [DateTime]"2021-12-10"
$ws.Cells.Item(12, 11).Value = "Castle Brite"
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 600
$ws.Cells.Item(12, 14).Value = 16000
$ws.Cells.Item(12, 15).Value = 16000
$ws.Cells.Item(12, 16).Value = 16000
$ws.Cells.Item(12, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(12, 18).Value = "Región del Maule"
$ws.Cells.Item(12, 19).Value = 889
$ws.Cells.Item(12, 20).Value = 18

# Row 13
$ws.Cells.Item(13, 4).Value = [DateTime]"2020-12-01"
$ws.Cells.Item(13, 11).Value = "Castle Brite"
$ws.Cells.Item(13, 12).Value = "Segunda"
$ws.Cells.Item(13, 13).Value = 600
$ws.Cells.Item(13, 14).Value = 16000
$ws.Cells.Item(13, 15).Value = 17000
$ws.Cells.Item(13, 16).Value = 16500
$ws.Cells.Item(13, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(13, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(13, 19).Value = 1100
$ws.Cells.Item(13, 20).Value = 15

# Row 14
$ws.Cells.Item(14, 4).Value = [DateTime]"2022-01-07"
$ws.Cells.Item(14, 11).Value = "Castle Brite"
$ws.Cells.Item(14, 12).Value = "Especial"
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 21000
$ws.Cells.Item(14, 15).Value = 21000
$ws.Cells.Item(14, 16).Value = 21000
$ws.Cells.Item(14, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(14, 18).Value = "Región Metropolitana"
$ws.Cells.Item(14, 19).Value = 1167
$ws.Cells.Item(14, 20).Value = 18

# Row 15
$ws.Cells.Item(15, 4).Value = [DateTime]"2022-01-07"
$ws.Cells.Item(15, 11).Value = "Castle Brite"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 18000
$ws.Cells.Item(15, 15).Value = 18000
$ws.Cells.Item(15, 16).Value = 18000
$ws.Cells.Item(15, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(15, 18).Value = "Región Metropolitana"
$ws.Cells.Item(15, 19).Value = 1000
$ws.Cells.Item(15, 20).Value = 18

# Row 16
$ws.Cells.Item(16, 4).Value = [DateTime]"2022-01-07"
$ws.Cells.Item(16, 11).Value = "Castle Brite"
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 16000
$ws.Cells.Item(16, 15).Value = 16000
$ws.Cells.Item(16, 16).Value = 16000
$ws.Cells.Item(16, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(16, 18).Value = "Región Metropolitana"
$ws.Cells.Item(16, 19).Value = 889
$ws.Cells.Item(16, 20).Value = 18

# Row 17
$ws.Cells.Item(17, 4).Value = [DateTime]"2021-12-03"
$ws.Cells.Item(17, 11).Value = "Castle Brite"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 350
$ws.Cells.Item(17, 14).Value = 24000
$ws.Cells.Item(17, 15).Value = 24000
$ws.Cells.Item(17, 16).Value = 24000
$ws.Cells.Item(17, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(17, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(17, 19).Value = 1333
$ws.Cells.Item(17, 20).Value = 18

# Row 18
$ws.Cells.Item(18, 4).Value = [DateTime]"2021-12-03"
$ws.Cells.Item(18, 11).Value = "Castle Brite"
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 350
$ws.Cells.Item(18, 14).Value = 20000
$ws.Cells.Item(18, 15).Value = 20000
$ws.Cells.Item(18, 16).Value = 20000
$ws.Cells.Item(18, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(18, 19).Value = 1111
$ws.Cells.Item(18, 20).Value = 18

# Row 19
$ws.Cells.Item(19, 4).Value = [DateTime]"2021-12-03"
$ws.Cells.Item(19, 11).Value = "Castle Brite"
$ws.Cells.Item(19, 12).Value = "Tercera"
$ws.Cells.Item(19, 13).Value = 350
$ws.Cells.Item(19, 14).Value = 17000
$ws.Cells.Item(19, 15).Value = 17000
$ws.Cells.Item(19, 16).Value = 17000
$ws.Cells.Item(19, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(19, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(19, 19).Value = 944
$ws.Cells.Item(19, 20).Value = 18

# Row 20
$ws.Cells.Item(20, 4).Value = [DateTime]"2020-11-27"
$ws.Cells.Item(20, 11).Value = "Castle Brite"
$ws.Cells.Item(20, 12).Value = "Tercera"
$ws.Cells.Item(20, 13).Value = 500
$ws.Cells.Item(20, 14).Value = 15000
$ws.Cells.Item(20, 15).Value = 16000
$ws.Cells.Item(20, 16).Value = 15500
$ws.Cells.Item(20, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 1033
$ws.Cells.Item(20, 20).Value = 15

# Row 21
$ws.Cells.Item(21, 4).Value = [DateTime]"2022-01-18"
$ws.Cells.Item(21, 11).Value = "Modesto"
$ws.Cells.Item(21, 12).Value = "Especial"
$ws.Cells.Item(21, 13).Value = 200
$ws.Cells.Item(21, 14).Value = 21000
$ws.Cells.Item(21, 15).Value = 21000
$ws.Cells.Item(21, 16).Value = 21000
$ws.Cells.Item(21, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(21, 18).Value = "Región Metropolitana"
$ws.Cells.Item(21, 19).Value = 1167
$ws.Cells.Item(21, 20).Value = 18

# Row 22
$ws.Cells.Item(22, 4).Value = [DateTime]"2022-01-18"
$ws.Cells.Item(22, 11).Value = "Modesto"
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 200
$ws.Cells.Item(22, 14).Value = 18000
$ws.Cells.Item(22, 15).Value = 18000
$ws.Cells.Item(22, 16).Value = 18000
$ws.Cells.Item(22, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(22, 18).Value = "Región Metropolitana"
$ws.Cells.Item(22, 19).Value = 1000
$ws.Cells.Item(22, 20).Value = 18

# Row 23
$ws.Cells.Item(23, 4).Value = [DateTime]"2022-01-18"
$ws.Cells.Item(23, 11).Value = "Modesto"
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 200
$ws.Cells.Item(23, 14).Value = 16000
$ws.Cells.Item(23, 15).Value = 16000
$ws.Cells.Item(23, 16).Value = 16000
$ws.Cells.Item(23, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(23, 18).Value = "Región Metropolitana"
$ws.Cells.Item(23, 19).Value = 889
$ws.Cells.Item(23, 20).Value = 18

# Row 24
$ws.Cells.Item(24, 4).Value = [DateTime]"2021-12-21"
$ws.Cells.Item(24, 11).Value = "Castle Brite"
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 200
$ws.Cells.Item(24, 14).Value = 20000
$ws.Cells.Item(24, 15).Value = 20000
$ws.Cells.Item(24, 16).Value = 20000
$ws.Cells.Item(24, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(24, 18).Value = "Región Metropolitana"
$ws.Cells.Item(24, 19).Value = 1111
$ws.Cells.Item(24, 20).Value = 18

# Row 25
$ws.Cells.Item(25, 4).Value = [DateTime]"2021-12-21"
$ws.Cells.Item(25, 11).Value = "Castle Brite"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 200
$ws.Cells.Item(25, 14).Value = 18000
$ws.Cells.Item(25, 15).Value = 18000
$ws.Cells.Item(25, 16).Value = 18000
$ws.Cells.Item(25, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(25, 18).Value = "Región Metropolitana"
$ws.Cells.Item(25, 19).Value = 1000
$ws.Cells.Item(25, 20).Value = 18

# Row 26
$ws.Cells.Item(26, 4).Value = [DateTime]"2021-12-21"
$ws.Cells.Item(26, 11).Value = "Castle Brite"
$ws.Cells.Item(26, 12).Value = "Segunda"
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 16000
$ws.Cells.Item(26, 15).Value = 16000
$ws.Cells.Item(26, 16).Value = 16000
$ws.Cells.Item(26, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(26, 18).Value = "Región Metropolitana"
$ws.Cells.Item(26, 19).Value = 889
$ws.Cells.Item(26, 20).Value = 18

# Row 27
$ws.Cells.Item(27, 4).Value = [DateTime]"2022-01-11"
$ws.Cells.Item(27, 11).Value = "Modesto"
$ws.Cells.Item(27, 12).Value = "Especial"
$ws.Cells.Item(27, 13).Value = 150
$ws.Cells.Item(27, 14).Value = 21000
$ws.Cells.Item(27, 15).Value = 21000
$ws.Cells.Item(27, 16).Value = 21000
$ws.Cells.Item(27, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(27, 18).Value = "Región Metropolitana"
$ws.Cells.Item(27, 19).Value = 1167
$ws.Cells.Item(27, 20).Value = 18

# Row 28
$ws.Cells.Item(28, 4).Value = [DateTime]"2022-01-11"
$ws.Cells.Item(28, 11).Value = "Modesto"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 150
$ws.Cells.Item(28, 14).Value = 18000
$ws.Cells.Item(28, 15).Value = 18000
$ws.Cells.Item(28, 16).Value = 18000
$ws.Cells.Item(28, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(28, 18).Value = "Región Metropolitana"
$ws.Cells.Item(28, 19).Value = 1000
$ws.Cells.Item(28, 20).Value = 18

# Row 29
$ws.Cells.Item(29, 4).Value = [DateTime]"2022-01-11"
$ws.Cells.Item(29, 11).Value = "Modesto"
$ws.Cells.Item(29, 12).Value = "Segunda"
$ws.Cells.Item(29, 13).Value = 150
$ws.Cells.Item(29, 14).Value = 16000
$ws.Cells.Item(29, 15).Value = 16000
$ws.Cells.Item(29, 16).Value = 16000
$ws.Cells.Item(29, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(29, 18).Value = "Región Metropolitana"
$ws.Cells.Item(29, 19).Value = 889
$ws.Cells.Item(29, 20).Value = 18

# Row 30
$ws.Cells.Item(30, 4).Value = [DateTime]"2021-11-30"
$ws.Cells.Item(30, 11).Value = "Castle Brite"
$ws.Cells.Item(30, 12).Value = "Segunda"
$ws.Cells.Item(30, 13).Value = 500
$ws.Cells.Item(30, 14).Value = 20000
$ws.Cells.Item(30, 15).Value = 21000
$ws.Cells.Item(30, 16).Value = 20500
$ws.Cells.Item(30, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(30, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(30, 19).Value = 1139
$ws.Cells.Item(30, 20).Value = 18

# Row 31
$ws.Cells.Item(31, 4).Value = [DateTime]"2021-12-23"
$ws.Cells.Item(31, 11).Value = "Castle Brite"
$ws.Cells.Item(31, 12).Value = "Especial"
$ws.Cells.Item(31, 13).Value = 250
$ws.Cells.Item(31, 14).Value = 20000
$ws.Cells.Item(31, 15).Value = 20000
$ws.Cells.Item(31, 16).Value = 20000
$ws.Cells.Item(31, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(31, 18).Value = "Región Metropolitana"
$ws.Cells.Item(31, 19).Value = 1111
$ws.Cells.Item(31, 20).Value = 18

# Row 32
$ws.Cells.Item(32, 4).Value = [DateTime]"2021-12-23"
$ws.Cells.Item(32, 11).Value = "Castle Brite"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 250
$ws.Cells.Item(32, 14).Value = 18000
$ws.Cells.Item(32, 15).Value = 18000
$ws.Cells.Item(32, 16).Value = 18000
$ws.Cells.Item(32, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(32, 18).Value = "Región Metropolitana"
$ws.Cells.Item(32, 19).Value = 1000
$ws.Cells.Item(32, 20).Value = 18

# Row 33
$ws.Cells.Item(33, 4).Value = [DateTime]"2021-12-23"
$ws.Cells.Item(33, 11).Value = "Castle Brite"
$ws.Cells.Item(33, 12).Value = "Segunda"
$ws.Cells.Item(33, 13).Value = 250
$ws.Cells.Item(33, 14).Value = 16000
$ws.Cells.Item(33, 15).Value = 16000
$ws.Cells.Item(33, 16).Value = 16000
$ws.Cells.Item(33, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(33, 18).Value = "Región Metropolitana"
$ws.Cells.Item(33, 19).Value = 889
$ws.Cells.Item(33, 20).Value = 18

# Row 34
$ws.Cells.Item(34, 4).Value = [DateTime]"2022-01-04"
$ws.Cells.Item(34, 11).Value = "Castle Brite"
$ws.Cells.Item(34, 12).Value = "Especial"
$ws.Cells.Item(34, 13).Value = 200
$ws.Cells.Item(34, 14).Value = 20000
$ws.Cells.Item(34, 15).Value = 20000
$ws.Cells.Item(34, 16).Value = 20000
$ws.Cells.Item(34, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(34, 18).Value = "Región Metropolitana"
$ws.Cells.Item(34, 19).Value = 1111
$ws.Cells.Item(34, 20).Value = 18

# Row 35
$ws.Cells.Item(35, 4).Value = [DateTime]"2022-01-04"
$ws.Cells.Item(35, 11).Value = "Castle Brite"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 200
$ws.Cells.Item(35, 14).Value = 18000
$ws.Cells.Item(35, 15).Value = 18000
$ws.Cells.Item(35, 16).Value = 18000
$ws.Cells.Item(35, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(35, 18).Value = "Región Metropolitana"
$ws.Cells.Item(35, 19).Value = 1000
$ws.Cells.Item(35, 20).Value = 18

# Row 36
$ws.Cells.Item(36, 4).Value = [DateTime]"2022-01-04"
$ws.Cells.Item(36, 11).Value = "Castle Brite"
$ws.Cells.Item(36, 12).Value = "Segunda"
$ws.Cells.Item(36, 13).Value = 200
$ws.Cells.Item(36, 14).Value = 16000
$ws.Cells.Item(36, 15).Value = 16000
$ws.Cells.Item(36, 16).Value = 16000
$ws.Cells.Item(36, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(36, 18).Value = "Región Metropolitana"
$ws.Cells.Item(36, 19).Value = 889
$ws.Cells.Item(36, 20).Value = 18

# Row 37
$ws.Cells.Item(37, 4).Value = [DateTime]"2020-12-22"
$ws.Cells.Item(37, 11).Value = "Castle Brite"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 350
$ws.Cells.Item(37, 14).Value = 16000
$ws.Cells.Item(37, 15).Value = 16000
$ws.Cells.Item(37, 16).Value = 16000
$ws.Cells.Item(37, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(37, 18).Value = "Región Metropolitana"
$ws.Cells.Item(37, 19).Value = 1067
$ws.Cells.Item(37, 20).Value = 15

# Row 38
$ws.Cells.Item(38, 4).Value = [DateTime]"2020-12-22"
$ws.Cells.Item(38, 11).Value = "Castle Brite"
$ws.Cells.Item(38, 12).Value = "Segunda"
$ws.Cells.Item(38, 13).Value = 300
$ws.Cells.Item(38, 14).Value = 13000
$ws.Cells.Item(38, 15).Value = 13000
$ws.Cells.Item(38, 16).Value = 13000
$ws.Cells.Item(38, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(38, 18).Value = "Región Metropolitana"
$ws.Cells.Item(38, 19).Value = 867
$ws.Cells.Item(38, 20).Value = 15

# Row 39
$ws.Cells.Item(39, 4).Value = [DateTime]"2021-12-14"
$ws.Cells.Item(39, 11).Value = "Castle Brite"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 600
$ws.Cells.Item(39, 14).Value = 18000
$ws.Cells.Item(39, 15).Value = 20000
$ws.Cells.Item(39, 16).Value = 19000
$ws.Cells.Item(39, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(39, 18).Value = "Región Metropolitana"
$ws.Cells.Item(39, 19).Value = 1056
$ws.Cells.Item(39, 20).Value = 18

# Row 40
$ws.Cells.Item(40, 4).Value = [DateTime]"2021-12-14"
$ws.Cells.Item(40, 11).Value = "Castle Brite"
$ws.Cells.Item(40, 12).Value = "Segunda"
$ws.Cells.Item(40, 13).Value = 300
$ws.Cells.Item(40, 14).Value = 16000
$ws.Cells.Item(40, 15).Value = 16000
$ws.Cells.Item(40, 16).Value = 16000
$ws.Cells.Item(40, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(40, 18).Value = "Región Metropolitana"
$ws.Cells.Item(40, 19).Value = 889
$ws.Cells.Item(40, 20).Value = 18

# Row 41
$ws.Cells.Item(41, 4).Value = [DateTime]"2021-12-07"
$ws.Cells.Item(41, 11).Value = "Castle Brite"
$ws.Cells.Item(41, 12).Value = "Primera"
$ws.Cells.Item(41, 13).Value = 500
$ws.Cells.Item(41, 14).Value = 20000
$ws.Cells.Item(41, 15).Value = 22000
$ws.Cells.Item(41, 16).Value = 21000
$ws.Cells.Item(41, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(41, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(41, 19).Value = 1167
$ws.Cells.Item(41, 20).Value = 18

# Row 42
$ws.Cells.Item(42, 4).Value = [DateTime]"2021-12-07"
$ws.Cells.Item(42, 11).Value = "Castle Brite"
$ws.Cells.Item(42, 12).Value = "Segunda"
$ws.Cells.Item(42, 13).Value = 250
$ws.Cells.Item(42, 14).Value = 17000
$ws.Cells.Item(42, 15).Value = 17000
$ws.Cells.Item(42, 16).Value = 17000
$ws.Cells.Item(42, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(42, 18).Value = "Región del Maule"
$ws.Cells.Item(42, 19).Value = 944
$ws.Cells.Item(42, 20).Value = 18
